# Fruta / hortaliza, semanal
# Insert two new weekly price rows (dated 45239) at the top of the
# "Palta" data block (right before the former row 1190), pushing the
# existing rows 1190:1221 down to 1192:1223.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 1190:1221 down by two rows.
$ws.Rows("1190:1191").Insert()

# New row 1190 - Palta Hass, "Primera" quality.
$ws.Cells.Item(1190, 1).Value  = 7
$ws.Cells.Item(1190, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(1190, 3).Value  = 'Ñuble'
$ws.Cells.Item(1190, 4).Value  = 45239
$ws.Cells.Item(1190, 5).Value  = 16
$ws.Cells.Item(1190, 6).Value  = 'Fruta'
$ws.Cells.Item(1190, 7).Value  = 100106
$ws.Cells.Item(1190, 8).Value  = 'Oleaginosos'
$ws.Cells.Item(1190, 9).Value  = 100106002
$ws.Cells.Item(1190, 10).Value = 'Palta'
$ws.Cells.Item(1190, 11).Value = 'Hass'
$ws.Cells.Item(1190, 12).Value = 'Primera'
$ws.Cells.Item(1190, 13).Value = 200
$ws.Cells.Item(1190, 14).Value = 3000
$ws.Cells.Item(1190, 15).Value = 3000
$ws.Cells.Item(1190, 16).Value = 3000
$ws.Cells.Item(1190, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1190, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1190, 19).Value = 3000
$ws.Cells.Item(1190, 20).Value = 1

# New row 1191 - Palta Hass, "Segunda" quality.
$ws.Cells.Item(1191, 1).Value  = 7
$ws.Cells.Item(1191, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(1191, 3).Value  = 'Ñuble'
$ws.Cells.Item(1191, 4).Value  = 45239
$ws.Cells.Item(1191, 5).Value  = 16
$ws.Cells.Item(1191, 6).Value  = 'Fruta'
$ws.Cells.Item(1191, 7).Value  = 100106
$ws.Cells.Item(1191, 8).Value  = 'Oleaginosos'
$ws.Cells.Item(1191, 9).Value  = 100106002
$ws.Cells.Item(1191, 10).Value = 'Palta'
$ws.Cells.Item(1191, 11).Value = 'Hass'
$ws.Cells.Item(1191, 12).Value = 'Segunda'
$ws.Cells.Item(1191, 13).Value = 150
$ws.Cells.Item(1191, 14).Value = 2800
$ws.Cells.Item(1191, 15).Value = 2800
$ws.Cells.Item(1191, 16).Value = 2800
$ws.Cells.Item(1191, 17).Value = '$/kilo (en caja de 17 kilos)'
$ws.Cells.Item(1191, 18).Value = 'Provincia de Quillota'
$ws.Cells.Item(1191, 19).Value = 2800
$ws.Cells.Item(1191, 20).Value = 1
